$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Macroferia Regional de Talca - Haba, Región del
# Maule origin, 20/10/2023) was inserted as row 49, pushing all subsequent
# rows (old 49..133) down by one (new 50..134).
$ws.Rows(49).Insert()

$ws.Cells.Item(49, 1).Value2 = 5
$ws.Cells.Item(49, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(49, 3).Value2 = "Maule"
$ws.Cells.Item(49, 4).Value2 = 45219
$ws.Cells.Item(49, 5).Value2 = 7
$ws.Cells.Item(49, 6).Value2 = 100112026
$ws.Cells.Item(49, 7).Value2 = "Haba"
$ws.Cells.Item(49, 8).Value2 = "Sin especificar"
$ws.Cells.Item(49, 9).Value2 = "Primera"
$ws.Cells.Item(49, 10).Value2 = 500
$ws.Cells.Item(49, 11).Value2 = 10000
$ws.Cells.Item(49, 12).Value2 = 10000
$ws.Cells.Item(49, 13).Value2 = 10000
$ws.Cells.Item(49, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value2 = "Región del Maule"
$ws.Cells.Item(49, 16).Value2 = 400
$ws.Cells.Item(49, 17).Value2 = 25
$ws.Cells.Item(49, 18).Value2 = "Hortaliza"
